# Lunggo_Config.xlsx edit:
#  - Add a new "airportFileName" config row (row 5), pushing the existing
#    "hotelLocationFileName" row down to row 6 (and everything below it
#    down by one row).
#  - Re-point hyperlinks that live on the shifted rows.
#  - Update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hyperlinks don't auto-follow a row insert in this engine, so clear them
# up front and recreate them afterwards at their new locations.
$ws.Hyperlinks.Delete()

# Insert a new blank row at row 6; existing row 6 (veritrans/Authorization)
# and everything below shifts down to row 7, etc.
$ws.Rows(6).Insert()

# Copy the formatting of row 4 (same visual style as old row 5) onto the
# newly-inserted row 6 so it matches the rest of the table.
$ws.Range("A4:H4").Copy() | Out-Null
$ws.Range("A6:H6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 5 becomes the new "airportFileName" entry.
$ws.Range("C5").Value = "airportFileName"
$ws.Range("E5").Value = "Airport.csv"
$ws.Range("F5").Value = "Airport.csv"
$ws.Range("G5").Value = "Airport.csv"

# Row 6 (the newly-inserted row) gets what used to be row 5's content.
$ws.Range("A6").Value = "*"
$ws.Range("B6").Value = "general"
$ws.Range("C6").Value = "hotelLocationFileName"
$ws.Range("D6").Formula = '="@@."&A6&"."&B6&"."&C6&"@@"'
$ws.Range("E6").Value = "HotelLocation.csv"
$ws.Range("F6").Value = "HotelLocation.csv"
$ws.Range("G6").Value = "HotelLocation.csv"

# Recreate the hyperlinks, shifted down by one row from row 6 onward.
$ws.Hyperlinks.Add($ws.Range("E8"), "https://api.sandbox.veritrans.co.id/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://api.sandbox.veritrans.co.id/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G8"), "https://api.sandbox.veritrans.co.id/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E13"), "http://travelmadezy.freshdesk.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G13"), "http://travelmadezy.freshdesk.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "http://travelmadezy.freshdesk.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E14"), "https://travelmadezy.zendesk.com/api/v2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E15"), "mailto:developer@travelmadezy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "mailto:developer@travelmadezy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G15"), "mailto:developer@travelmadezy.com") | Out-Null

# Update selection/scroll state to match.
$ws.Range("A1").Select()
$ws.Range("D12").Select()
